{"js": "// Apply the benchmark-table update described by the commit:\n// - first 4 rows get new summary values (\"0M\",\"0M\",\"0M\",\"23\")\n// - the next several rows of percentages get updated\n// - the last 3 rows (which held tab-separated per-iteration detail\n//   values) get collapsed down to the single leading value.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Row index (0-based) -> new single-cell text value.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"23\",\n  4: \"0.00002\",\n  5: \"0.00005\",\n  6: \"0.00003\",\n  7: \"0.00001\",\n  8: \"0.00003\",\n  9: \"0.00003\",\n  10: \"0.00004\",\n  11: \"0.00076\",\n  43: \"100\",\n  44: \"0\",\n  45: \"27\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(parseInt(rowIndex, 10), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-table update described by the commit:\n# - first 4 rows get new summary values (\"0M\",\"0M\",\"0M\",\"23\")\n# - the next several rows of percentages get updated\n# - the last 3 rows (which held tab-separated per-iteration detail\n#   values) get collapsed down to the single leading value.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Map of 1-based row index -> new single-cell text value.\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"23\"\n    5  = \"0.00002\"\n    6  = \"0.00005\"\n    7  = \"0.00003\"\n    8  = \"0.00001\"\n    9  = \"0.00003\"\n    10 = \"0.00003\"\n    11 = \"0.00004\"\n    12 = \"0.00076\"\n    44 = \"100\"\n    45 = \"0\"\n    46 = \"27\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $tbl.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
